$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date strings in column B from the 03/08/2023-style values to
# the corresponding 05/08/2023-style values, preserving each row's
# existing text format (leading "'" keeps the cell text, not a date
# serial, so the original cell style is not disturbed).
for ($r = 1; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value()
    if ($v -eq "03/08/2023") {
        $cell.Value = "'05/08/2023"
    } elseif ($v -eq "2023/08/03") {
        $cell.Value = "'2023/08/05"
    } elseif ($v -eq "08/03/2023") {
        $cell.Value = "'08/05/2023"
    }
}

# Move the active selection from B17 to C13
$ws.Range("C13").Select()
